# Egypt Division 1 - odds feed refresh (11-04-2024 00:31)
# Applies: (1) three row-groups whose match records were
# reordered by the source feed (values only; id/date columns
# already matched), (2) closing-line refresh on two upcoming
# fixtures, and (3) one brand-new upcoming fixture appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Egypt Division 1")

function Set-Cell($addr, $val) {
    $ws.Range($addr).Value = $val
}

# --- Re-ordered match rows (same id/date columns, B:AC payload swapped) ---
# Row 39
Set-Cell "B39" 6853141
Set-Cell "C39" "Egypt Division 1"
Set-Cell "D39" "Egypt Division 1"
Set-Cell "F39" "Aswan FC"
Set-Cell "G39" "El Daklyeh"
Set-Cell "H39" 2
Set-Cell "I39" 4
Set-Cell "J39" "A"
Set-Cell "K39" 2
Set-Cell "L39" 2.75
Set-Cell "M39" 4
Set-Cell "N39" 2.4
Set-Cell "O39" 2.75
Set-Cell "P39" 3
Set-Cell "Q39" -0.25
Set-Cell "R39" 2
Set-Cell "S39" 1.8
Set-Cell "T39" 2.25
Set-Cell "U39" 1.9
Set-Cell "V39" 1.9
Set-Cell "W39" -1
Set-Cell "X39" -1
Set-Cell "Y39" 2
Set-Cell "Z39" -1
Set-Cell "AA39" 0.8
Set-Cell "AB39" 0.8999999999999999
Set-Cell "AC39" -1

# Row 40
Set-Cell "B40" 6853132
Set-Cell "C40" "Egypt Division 1"
Set-Cell "D40" "Egypt Division 1"
Set-Cell "F40" "Haras El Hedoud"
Set-Cell "G40" "Ghazl El Mahallah"
Set-Cell "H40" 2
Set-Cell "I40" 1
Set-Cell "J40" "H"
Set-Cell "K40" 5.5
Set-Cell "L40" 3.25
Set-Cell "M40" 1.6
Set-Cell "N40" 4.5
Set-Cell "O40" 3.4
Set-Cell "P40" 1.666
Set-Cell "Q40" 0.75
Set-Cell "R40" 1.825
Set-Cell "S40" 1.975
Set-Cell "T40" 2.25
Set-Cell "U40" 1.8
Set-Cell "V40" 2
Set-Cell "W40" 3.5
Set-Cell "X40" -1
Set-Cell "Y40" -1
Set-Cell "Z40" 0.825
Set-Cell "AA40" -1
Set-Cell "AB40" 0.8
Set-Cell "AC40" -1

# Row 41
Set-Cell "B41" 6853133
Set-Cell "C41" "Egypt Division 1"
Set-Cell "D41" "Egypt Division 1"
Set-Cell "F41" "Pyramids FC"
Set-Cell "G41" "Talaea El Geish"
Set-Cell "H41" 4
Set-Cell "I41" 2
Set-Cell "J41" "H"
Set-Cell "K41" 1.727
Set-Cell "L41" 3.25
Set-Cell "M41" 4.5
Set-Cell "N41" 1.6
Set-Cell "O41" 3.3
Set-Cell "P41" 5.5
Set-Cell "Q41" -1
Set-Cell "R41" 1.975
Set-Cell "S41" 1.825
Set-Cell "T41" 2.5
Set-Cell "U41" 1.9
Set-Cell "V41" 1.9
Set-Cell "W41" 0.6000000000000001
Set-Cell "X41" -1
Set-Cell "Y41" -1
Set-Cell "Z41" 0.9750000000000001
Set-Cell "AA41" -1
Set-Cell "AB41" 0.8999999999999999
Set-Cell "AC41" -1

# Row 54
Set-Cell "B54" 7208758
Set-Cell "C54" "Egypt Division 1"
Set-Cell "D54" "Egypt Division 1"
Set-Cell "F54" "Smouha"
Set-Cell "G54" "El Gounah"
Set-Cell "H54" 1
Set-Cell "I54" 1
Set-Cell "J54" "D"
Set-Cell "K54" 2.5
Set-Cell "L54" 2.8
Set-Cell "M54" 2.8
Set-Cell "N54" 2.5
Set-Cell "O54" 2.625
Set-Cell "P54" 3
Set-Cell "Q54" 0
Set-Cell "R54" 1.725
Set-Cell "S54" 2.075
Set-Cell "T54" 2
Set-Cell "U54" 1.775
Set-Cell "V54" 2.025
Set-Cell "W54" -1
Set-Cell "X54" 1.625
Set-Cell "Y54" -1
Set-Cell "Z54" 0
Set-Cell "AA54" -0
Set-Cell "AB54" 0
Set-Cell "AC54" -0

# Row 55
Set-Cell "B55" 7210310
Set-Cell "C55" "Egypt Division 1"
Set-Cell "D55" "Egypt Division 1"
Set-Cell "F55" "Al Moqawloon Al Arab"
Set-Cell "G55" "Baladiyet El Mahallah"
Set-Cell "H55" 0
Set-Cell "I55" 2
Set-Cell "J55" "A"
Set-Cell "K55" 2.5
Set-Cell "L55" 2.7
Set-Cell "M55" 2.9
Set-Cell "N55" 2.25
Set-Cell "O55" 2.7
Set-Cell "P55" 3.5
Set-Cell "Q55" -0.25
Set-Cell "R55" 1.925
Set-Cell "S55" 1.875
Set-Cell "T55" 1.75
Set-Cell "U55" 1.925
Set-Cell "V55" 1.875
Set-Cell "W55" -1
Set-Cell "X55" -1
Set-Cell "Y55" 2.5
Set-Cell "Z55" -1
Set-Cell "AA55" 0.875
Set-Cell "AB55" 0.4625
Set-Cell "AC55" -0.5

# Row 57
Set-Cell "B57" 7208367
Set-Cell "C57" "Egypt Division 1"
Set-Cell "D57" "Egypt Division 1"
Set-Cell "F57" "Pyramids FC"
Set-Cell "G57" "El Zamalek"
Set-Cell "H57" 2
Set-Cell "I57" 2
Set-Cell "J57" "D"
Set-Cell "K57" 1.909
Set-Cell "L57" 3.2
Set-Cell "M57" 3.8
Set-Cell "N57" 2.05
Set-Cell "O57" 3.1
Set-Cell "P57" 3.5
Set-Cell "Q57" -0.5
Set-Cell "R57" 2.025
Set-Cell "S57" 1.775
Set-Cell "T57" 2.5
Set-Cell "U57" 1.975
Set-Cell "V57" 1.825
Set-Cell "W57" -1
Set-Cell "X57" 2.1
Set-Cell "Y57" -1
Set-Cell "Z57" -1
Set-Cell "AA57" 0.7749999999999999
Set-Cell "AB57" 0.9750000000000001
Set-Cell "AC57" -1

# Row 58
Set-Cell "B58" 7208756
Set-Cell "C58" "Egypt Division 1"
Set-Cell "D58" "Egypt Division 1"
Set-Cell "F58" "Al Ittihad Al Sakandary"
Set-Cell "G58" "Ceramica Cleopatra"
Set-Cell "H58" 1
Set-Cell "I58" 0
Set-Cell "J58" "H"
Set-Cell "K58" 2.8
Set-Cell "L58" 2.9
Set-Cell "M58" 2.55
Set-Cell "N58" 2.8
Set-Cell "O58" 2.9
Set-Cell "P58" 2.55
Set-Cell "Q58" 0
Set-Cell "R58" 1.975
Set-Cell "S58" 1.825
Set-Cell "T58" 2.25
Set-Cell "U58" 2
Set-Cell "V58" 1.8
Set-Cell "W58" 1.8
Set-Cell "X58" -1
Set-Cell "Y58" -1
Set-Cell "Z58" 0.9750000000000001
Set-Cell "AA58" -1
Set-Cell "AB58" -1
Set-Cell "AC58" 0.8

# Row 119
Set-Cell "B119" 7217680
Set-Cell "C119" "Egypt Division 1"
Set-Cell "D119" "Egypt Division 1"
Set-Cell "F119" "El Masry"
Set-Cell "G119" "El Zamalek"
Set-Cell "H119" 1
Set-Cell "I119" 0
Set-Cell "J119" "H"
Set-Cell "K119" 2.3
Set-Cell "L119" 3
Set-Cell "M119" 2.875
Set-Cell "N119" 2.875
Set-Cell "O119" 3.2
Set-Cell "P119" 2.375
Set-Cell "Q119" 0.25
Set-Cell "R119" 1.75
Set-Cell "S119" 2.05
Set-Cell "T119" 2.5
Set-Cell "U119" 2.025
Set-Cell "V119" 1.775
Set-Cell "W119" 1.875
Set-Cell "X119" -1
Set-Cell "Y119" -1
Set-Cell "Z119" 0.75
Set-Cell "AA119" -1
Set-Cell "AB119" -1
Set-Cell "AC119" 0.7749999999999999

# Row 120
Set-Cell "B120" 7217677
Set-Cell "C120" "Egypt Division 1"
Set-Cell "D120" "Egypt Division 1"
Set-Cell "F120" "Pyramids FC"
Set-Cell "G120" "Baladiyet El Mahallah"
Set-Cell "H120" 2
Set-Cell "I120" 2
Set-Cell "J120" "D"
Set-Cell "K120" 1.25
Set-Cell "L120" 5
Set-Cell "M120" 9
Set-Cell "N120" 1.166
Set-Cell "O120" 7
Set-Cell "P120" 11
Set-Cell "Q120" -2
Set-Cell "R120" 1.95
Set-Cell "S120" 1.85
Set-Cell "T120" 3
Set-Cell "U120" 1.825
Set-Cell "V120" 1.975
Set-Cell "W120" -1
Set-Cell "X120" 6
Set-Cell "Y120" -1
Set-Cell "Z120" -1
Set-Cell "AA120" 0.8500000000000001
Set-Cell "AB120" 0.825
Set-Cell "AC120" -1

# Row 131
Set-Cell "B131" 7217692
Set-Cell "C131" "Egypt Division 1"
Set-Cell "D131" "Egypt Division 1"
Set-Cell "F131" "Talaea El Geish"
Set-Cell "G131" "ZED FC"
Set-Cell "H131" 1
Set-Cell "I131" 0
Set-Cell "J131" "H"
Set-Cell "K131" 2.9
Set-Cell "L131" 2.9
Set-Cell "M131" 2.4
Set-Cell "N131" 2.8
Set-Cell "O131" 2.875
Set-Cell "P131" 2.5
Set-Cell "Q131" 0
Set-Cell "R131" 1.975
Set-Cell "S131" 1.825
Set-Cell "T131" 2
Set-Cell "U131" 2.025
Set-Cell "V131" 1.775
Set-Cell "W131" 1.8
Set-Cell "X131" -1
Set-Cell "Y131" -1
Set-Cell "Z131" 0.9750000000000001
Set-Cell "AA131" -1
Set-Cell "AB131" -1
Set-Cell "AC131" 0.7749999999999999

# Row 132
Set-Cell "B132" 7217691
Set-Cell "C132" "Egypt Division 1"
Set-Cell "D132" "Egypt Division 1"
Set-Cell "F132" "Enppi"
Set-Cell "G132" "Baladiyet El Mahallah"
Set-Cell "H132" 3
Set-Cell "I132" 0
Set-Cell "J132" "H"
Set-Cell "K132" 2.1
Set-Cell "L132" 2.8
Set-Cell "M132" 3.6
Set-Cell "N132" 1.85
Set-Cell "O132" 3
Set-Cell "P132" 4.2
Set-Cell "Q132" -0.5
Set-Cell "R132" 1.9
Set-Cell "S132" 1.9
Set-Cell "T132" 2.25
Set-Cell "U132" 1.95
Set-Cell "V132" 1.85
Set-Cell "W132" 0.8500000000000001
Set-Cell "X132" -1
Set-Cell "Y132" -1
Set-Cell "Z132" 0.8999999999999999
Set-Cell "AA132" -1
Set-Cell "AB132" 0.95
Set-Cell "AC132" -1

# Row 137
Set-Cell "B137" 7217701
Set-Cell "C137" "Egypt Division 1"
Set-Cell "D137" "Egypt Division 1"
Set-Cell "F137" "Al Ittihad Al Sakandary"
Set-Cell "G137" "Talaea El Geish"
Set-Cell "H137" 2
Set-Cell "I137" 2
Set-Cell "J137" "D"
Set-Cell "K137" 2.375
Set-Cell "L137" 2.9
Set-Cell "M137" 3
Set-Cell "N137" 2.6
Set-Cell "O137" 3
Set-Cell "P137" 2.8
Set-Cell "Q137" 0
Set-Cell "R137" 1.825
Set-Cell "S137" 1.975
Set-Cell "T137" 2
Set-Cell "U137" 1.925
Set-Cell "V137" 1.875
Set-Cell "W137" -1
Set-Cell "X137" 2
Set-Cell "Y137" -1
Set-Cell "Z137" 0
Set-Cell "AA137" -0
Set-Cell "AB137" 0.925
Set-Cell "AC137" -1

# Row 138
Set-Cell "B138" 7217702
Set-Cell "C138" "Egypt Division 1"
Set-Cell "D138" "Egypt Division 1"
Set-Cell "F138" "Al Moqawloon Al Arab"
Set-Cell "G138" "Enppi"
Set-Cell "H138" 0
Set-Cell "I138" 1
Set-Cell "J138" "A"
Set-Cell "K138" 2.9
Set-Cell "L138" 3
Set-Cell "M138" 2.375
Set-Cell "N138" 2.75
Set-Cell "O138" 2.8
Set-Cell "P138" 2.5
Set-Cell "Q138" 0
Set-Cell "R138" 2.025
Set-Cell "S138" 1.775
Set-Cell "T138" 1.75
Set-Cell "U138" 1.775
Set-Cell "V138" 2.025
Set-Cell "W138" -1
Set-Cell "X138" -1
Set-Cell "Y138" 1.5
Set-Cell "Z138" -1
Set-Cell "AA138" 0.7749999999999999
Set-Cell "AB138" -1
Set-Cell "AC138" 1.025

# Row 169
Set-Cell "B169" 7217740
Set-Cell "C169" "Egypt Division 1"
Set-Cell "D169" "Egypt Division 1"
Set-Cell "F169" "Smouha"
Set-Cell "G169" "Ceramica Cleopatra"
Set-Cell "H169" 1
Set-Cell "I169" 0
Set-Cell "J169" "H"
Set-Cell "K169" 2.9
Set-Cell "L169" 2.8
Set-Cell "M169" 2.4
Set-Cell "N169" 3
Set-Cell "O169" 2.9
Set-Cell "P169" 2.25
Set-Cell "Q169" 0.25
Set-Cell "R169" 1.8
Set-Cell "S169" 2
Set-Cell "T169" 2.25
Set-Cell "U169" 1.825
Set-Cell "V169" 1.975
Set-Cell "W169" 2
Set-Cell "X169" -1
Set-Cell "Y169" -1
Set-Cell "Z169" 0.8
Set-Cell "AA169" -1
Set-Cell "AB169" -1
Set-Cell "AC169" 0.9750000000000001

# Row 170
Set-Cell "B170" 7217742
Set-Cell "C170" "Egypt Division 1"
Set-Cell "D170" "Egypt Division 1"
Set-Cell "F170" "El Daklyeh"
Set-Cell "G170" "Pyramids FC"
Set-Cell "H170" 0
Set-Cell "I170" 2
Set-Cell "J170" "A"
Set-Cell "K170" 7.5
Set-Cell "L170" 3.75
Set-Cell "M170" 1.4
Set-Cell "N170" 8.5
Set-Cell "O170" 4
Set-Cell "P170" 1.363
Set-Cell "Q170" 1.25
Set-Cell "R170" 1.875
Set-Cell "S170" 1.925
Set-Cell "T170" 2.25
Set-Cell "U170" 1.9
Set-Cell "V170" 1.9
Set-Cell "W170" -1
Set-Cell "X170" -1
Set-Cell "Y170" 0.363
Set-Cell "Z170" -1
Set-Cell "AA170" 0.925
Set-Cell "AB170" -0.5
Set-Cell "AC170" 0.45

# Row 175
Set-Cell "B175" 7881846
Set-Cell "C175" "Egypt Division 1"
Set-Cell "D175" "Egypt Division 1"
Set-Cell "F175" "Ceramica Cleopatra"
Set-Cell "G175" "ZED FC"
Set-Cell "H175" 1
Set-Cell "I175" 1
Set-Cell "J175" "D"
Set-Cell "K175" 1.909
Set-Cell "L175" 3.2
Set-Cell "M175" 3.8
Set-Cell "N175" 1.8
Set-Cell "O175" 3.4
Set-Cell "P175" 4
Set-Cell "Q175" -0.5
Set-Cell "R175" 1.775
Set-Cell "S175" 2.025
Set-Cell "T175" 2.5
Set-Cell "U175" 1.975
Set-Cell "V175" 1.825
Set-Cell "W175" -1
Set-Cell "X175" 2.4
Set-Cell "Y175" -1
Set-Cell "Z175" -1
Set-Cell "AA175" 1.025
Set-Cell "AB175" -1
Set-Cell "AC175" 0.825

# Row 176
Set-Cell "B176" 7217744
Set-Cell "C176" "Egypt Division 1"
Set-Cell "D176" "Egypt Division 1"
Set-Cell "F176" "El Gounah"
Set-Cell "G176" "Al Moqawloon Al Arab"
Set-Cell "H176" 1
Set-Cell "I176" 2
Set-Cell "J176" "A"
Set-Cell "K176" 2.45
Set-Cell "L176" 2.9
Set-Cell "M176" 2.9
Set-Cell "N176" 2.45
Set-Cell "O176" 2.875
Set-Cell "P176" 2.9
Set-Cell "Q176" 0
Set-Cell "R176" 1.75
Set-Cell "S176" 2.05
Set-Cell "T176" 2
Set-Cell "U176" 1.9
Set-Cell "V176" 1.9
Set-Cell "W176" -1
Set-Cell "X176" -1
Set-Cell "Y176" 1.9
Set-Cell "Z176" -1
Set-Cell "AA176" 1.05
Set-Cell "AB176" 0.8999999999999999
Set-Cell "AC176" -1

# Row 177
Set-Cell "B177" 7217745
Set-Cell "C177" "Egypt Division 1"
Set-Cell "D177" "Egypt Division 1"
Set-Cell "F177" "Ismaily SC"
Set-Cell "G177" "Smouha"
Set-Cell "H177" 2
Set-Cell "I177" 1
Set-Cell "J177" "H"
Set-Cell "K177" 2.7
Set-Cell "L177" 2.9
Set-Cell "M177" 2.6
Set-Cell "N177" 2.55
Set-Cell "O177" 2.9
Set-Cell "P177" 2.75
Set-Cell "Q177" 0
Set-Cell "R177" 1.825
Set-Cell "S177" 1.975
Set-Cell "T177" 2
Set-Cell "U177" 1.825
Set-Cell "V177" 1.975
Set-Cell "W177" 1.55
Set-Cell "X177" -1
Set-Cell "Y177" -1
Set-Cell "Z177" 0.825
Set-Cell "AA177" -1
Set-Cell "AB177" 0.825
Set-Cell "AC177" -1

# Row 178
Set-Cell "B178" 8031210
Set-Cell "C178" "Egypt Division 1"
Set-Cell "D178" "Egypt Division 1"
Set-Cell "F178" "Enppi"
Set-Cell "G178" "Pharco FC"
Set-Cell "H178" 1
Set-Cell "I178" 1
Set-Cell "J178" "D"
Set-Cell "K178" 2.1
Set-Cell "L178" 2.8
Set-Cell "M178" 3.8
Set-Cell "N178" 2.45
Set-Cell "O178" 2.8
Set-Cell "P178" 3
Set-Cell "Q178" -0.25
Set-Cell "R178" 2.075
Set-Cell "S178" 1.725
Set-Cell "T178" 1.75
Set-Cell "U178" 1.775
Set-Cell "V178" 2.025
Set-Cell "W178" -1
Set-Cell "X178" 1.8
Set-Cell "Y178" -1
Set-Cell "Z178" -0.5
Set-Cell "AA178" 0.3625
Set-Cell "AB178" 0.3875
Set-Cell "AC178" -0.5

# Row 179
Set-Cell "B179" 8031212
Set-Cell "C179" "Egypt Division 1"
Set-Cell "D179" "Egypt Division 1"
Set-Cell "F179" "National Bank"
Set-Cell "G179" "El Daklyeh"
Set-Cell "H179" 1
Set-Cell "I179" 0
Set-Cell "J179" "H"
Set-Cell "K179" 1.6
Set-Cell "L179" 3.5
Set-Cell "M179" 5.75
Set-Cell "N179" 1.7
Set-Cell "O179" 3.3
Set-Cell "P179" 5
Set-Cell "Q179" -0.75
Set-Cell "R179" 1.975
Set-Cell "S179" 1.825
Set-Cell "T179" 2.25
Set-Cell "U179" 1.875
Set-Cell "V179" 1.925
Set-Cell "W179" 0.7
Set-Cell "X179" -1
Set-Cell "Y179" -1
Set-Cell "Z179" 0.4875
Set-Cell "AA179" -0.5
Set-Cell "AB179" -1
Set-Cell "AC179" 0.925

# Row 180
Set-Cell "B180" 8031211
Set-Cell "C180" "Egypt Division 1"
Set-Cell "D180" "Egypt Division 1"
Set-Cell "F180" "El Masry"
Set-Cell "G180" "Al Ittihad Al Sakandary"
Set-Cell "H180" 2
Set-Cell "I180" 3
Set-Cell "J180" "A"
Set-Cell "K180" 1.95
Set-Cell "L180" 3
Set-Cell "M180" 4
Set-Cell "N180" 2
Set-Cell "O180" 3
Set-Cell "P180" 3.8
Set-Cell "Q180" -0.5
Set-Cell "R180" 2
Set-Cell "S180" 1.8
Set-Cell "T180" 2.25
Set-Cell "U180" 1.825
Set-Cell "V180" 1.975
Set-Cell "W180" -1
Set-Cell "X180" -1
Set-Cell "Y180" 2.8
Set-Cell "Z180" -1
Set-Cell "AA180" 0.8
Set-Cell "AB180" 0.825
Set-Cell "AC180" -1

# --- Closing-line refresh on upcoming fixtures (row 184 & 185) ---
Set-Cell "N184" 2.3
Set-Cell "O184" 2.9
Set-Cell "P184" 2.9
Set-Cell "R184" 2.05
Set-Cell "S184" 1.75

Set-Cell "N185" 1.85
Set-Cell "O185" 3.2
Set-Cell "P185" 3.75
Set-Cell "Q185" -0.5
Set-Cell "R185" 1.9
Set-Cell "S185" 1.9
Set-Cell "U185" 1.975
Set-Cell "V185" 1.825

# --- New upcoming fixture appended as row 186 ---
Set-Cell "A186" 184
Set-Cell "B186" 7880128
Set-Cell "C186" "Egypt Division 1"
Set-Cell "D186" "Egypt Division 1"
Set-Cell "E186" 45393.58333333334
Set-Cell "F186" "ZED FC"
Set-Cell "G186" "Al Ahly Cairo"
Set-Cell "K186" 5.5
Set-Cell "L186" 3.6
Set-Cell "M186" 1.533
Set-Cell "N186" 6.5
Set-Cell "O186" 3.6
Set-Cell "P186" 1.5
Set-Cell "Q186" 1
Set-Cell "R186" 1.95
Set-Cell "S186" 1.85
Set-Cell "T186" 2.25
Set-Cell "U186" 1.8
Set-Cell "V186" 2
Set-Cell "W186" 0
Set-Cell "X186" 0
Set-Cell "Y186" 0
Set-Cell "Z186" 0
Set-Cell "AA186" 0

# Row 186 needs the same formatting as the other data rows:
#   column A -> bold, centered, thin-bordered (style used by all A/id cells)
#   column E -> custom date/time number format (style used by all E/date cells)
$ws.Range("A185").Copy() | Out-Null
$ws.Range("A186").PasteSpecial(-4122) | Out-Null
$ws.Range("E185").Copy() | Out-Null
$ws.Range("E186").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
